$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4401.5454  # was 4272.206
$ws.Range("I11").Value = 4401.5454  # was 4272.206
$ws.Range("K11").Value = 4401.5454  # was 4272.206
$ws.Range("M11").Value = -4261.5454  # was -4132.206
$ws.Range("H74").Value = 8308.637000000001  # was 9399.700000000001
$ws.Range("I74").Value = 7342.143  # was 8799.4
$ws.Range("K74").Value = 7342.143  # was 8799.4
$ws.Range("M74").Value = -6406.143  # was -7863.4
$ws.Range("H77").Value = 8308.637000000001  # was 9399.700000000001
$ws.Range("I77").Value = 7342.143  # was 8799.4
$ws.Range("K77").Value = 36710.715  # was 43997
$ws.Range("M77").Value = -32030.715  # was -39317
$ws.Range("H82").Value = 262  # was 156.5
$ws.Range("I82").Value = 262  # was 156.5
$ws.Range("K82").Value = 786  # was 469.5
$ws.Range("M82").Value = -380  # was -63.5
$ws.Range("H85").Value = 262  # was 156.5
$ws.Range("I85").Value = 262  # was 156.5
$ws.Range("K85").Value = 786  # was 469.5
$ws.Range("M85").Value = 618  # was 934.5
$ws.Range("H86").Value = 4000  # was 4222
$ws.Range("J86").Value = 0  # was 4444
$ws.Range("L86").Value = 0  # was 4444
$ws.Range("N86").ClearContents()  # was -6690
$ws.Range("H88").Value = 6490.4  # was 6519.9
$ws.Range("I88").Value = 9350  # was 18500
$ws.Range("J88").Value = 5775.5  # was 5188.778
$ws.Range("K88").Value = 9350  # was 18500
$ws.Range("L88").Value = 5775.5  # was 5188.778
$ws.Range("M88").Value = -8944  # was -18094
$ws.Range("N88").Value = -6587.5  # was -6000.778
$ws.Range("H89").Value = 4000  # was 4222
$ws.Range("J89").Value = 0  # was 4444
$ws.Range("L89").Value = 0  # was 22220
$ws.Range("N89").ClearContents()  # was -33452
$ws.Range("H91").Value = 6490.4  # was 6519.9
$ws.Range("I91").Value = 9350  # was 18500
$ws.Range("J91").Value = 5775.5  # was 5188.778
$ws.Range("K91").Value = 9350  # was 18500
$ws.Range("L91").Value = 5775.5  # was 5188.778
$ws.Range("M91").Value = -7946  # was -17096
$ws.Range("N91").Value = -8583.5  # was -7996.778
$ws.Range("H127").Value = 1545.5264  # was 1543.8235
$ws.Range("I127").Value = 933.4545000000001  # was 944.8
$ws.Range("J127").Value = 2387.125  # was 2399.5715
$ws.Range("K127").Value = 2800.3635  # was 2834.4
$ws.Range("L127").Value = 7161.375  # was 7198.7145
$ws.Range("M127").Value = 2159.6365  # was 2125.6
$ws.Range("N127").Value = -17081.375  # was -17118.7145
$ws.Range("H135").Value = 33343700  # was 10004081
$ws.Range("I135").Value = 100000000  # was 14286516
$ws.Range("J135").Value = 15550  # was 11733.333
$ws.Range("K135").Value = 900000000  # was 128578644
$ws.Range("L135").Value = 139950  # was 105599.997
$ws.Range("M135").Value = -899997465  # was -128576109
$ws.Range("N135").Value = -145020  # was -110669.997

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 24850  # was 35212
$ws.Range("I55").Value = 0  # was 43048
$ws.Range("J55").Value = 24850  # was 32600
$ws.Range("K55").Value = 0  # was 43048
$ws.Range("L55").Value = 24850  # was 32600
$ws.Range("M55").ClearContents()  # was -42733
$ws.Range("N55").Value = -25480  # was -33230
$ws.Range("H63").Value = 9009.954  # was 7724.76
$ws.Range("I63").Value = 1902.5  # was 1367.5
$ws.Range("J63").Value = 9720.700000000001  # was 9732.315000000001
$ws.Range("K63").Value = 1902.5  # was 1367.5
$ws.Range("L63").Value = 9720.700000000001  # was 9732.315000000001
$ws.Range("M63").Value = -1216.5  # was -681.5
$ws.Range("N63").Value = -11092.7  # was -11104.315
$ws.Range("H66").Value = 9009.954  # was 7724.76
$ws.Range("I66").Value = 1902.5  # was 1367.5
$ws.Range("J66").Value = 9720.700000000001  # was 9732.315000000001
$ws.Range("K66").Value = 9512.5  # was 6837.5
$ws.Range("L66").Value = 48603.5  # was 48661.575
$ws.Range("M66").Value = -6080.5  # was -3405.5
$ws.Range("N66").Value = -55467.5  # was -55525.575
$ws.Range("H74").Value = 3988.577  # was 4343.7827
$ws.Range("I74").Value = 1744.3636  # was 1932.7778
$ws.Range("J74").Value = 5634.3335  # was 5893.7144
$ws.Range("K74").Value = 1744.3636  # was 1932.7778
$ws.Range("L74").Value = 5634.3335  # was 5893.7144
$ws.Range("M74").Value = -870.3635999999999  # was -1058.7778
$ws.Range("N74").Value = -7382.3335  # was -7641.7144
$ws.Range("H77").Value = 3988.577  # was 4343.7827
$ws.Range("I77").Value = 1744.3636  # was 1932.7778
$ws.Range("J77").Value = 5634.3335  # was 5893.7144
$ws.Range("K77").Value = 8721.817999999999  # was 9663.889000000001
$ws.Range("L77").Value = 28171.6675  # was 29468.572
$ws.Range("M77").Value = -4353.817999999999  # was -5295.889000000001
$ws.Range("N77").Value = -36907.6675  # was -38204.572
$ws.Range("H80").Value = 110072  # was 86702.664
$ws.Range("I80").Value = 0  # was 20000
$ws.Range("J80").Value = 110072  # was 120054
$ws.Range("K80").Value = 0  # was 20000
$ws.Range("L80").Value = 110072  # was 120054
$ws.Range("M80").ClearContents()  # was -19002
$ws.Range("N80").Value = -112068  # was -122050
$ws.Range("H83").Value = 110072  # was 86702.664
$ws.Range("I83").Value = 0  # was 20000
$ws.Range("J83").Value = 110072  # was 120054
$ws.Range("K83").Value = 0  # was 60000
$ws.Range("L83").Value = 330216  # was 360162
$ws.Range("M83").ClearContents()  # was -55008
$ws.Range("N83").Value = -340200  # was -370146
$ws.Range("H122").Value = 1609.5  # was 1577.1
$ws.Range("I122").Value = 1410.8572  # was 1419
$ws.Range("K122").Value = 4232.571599999999  # was 4257
$ws.Range("M122").Value = -1782.571599999999  # was -1807

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2032.5  # was 2546.25
$ws.Range("I99").Value = 2032.5  # was 2546.25
$ws.Range("K99").Value = 2032.5  # was 2546.25
$ws.Range("M99").Value = -534.5  # was -1048.25
$ws.Range("H107").Value = 2306.625  # was 2439.4285
$ws.Range("I107").Value = 2034.0869  # was 2132.65
$ws.Range("K107").Value = 2034.0869  # was 2132.65
$ws.Range("M107").Value = -114.0869  # was -212.6500000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 1971.6666  # was 500
$ws.Range("I21").Value = 0  # was 500
$ws.Range("J21").Value = 1971.6666  # was 0
$ws.Range("K21").Value = 0  # was 500
$ws.Range("L21").Value = 1971.6666  # was 0
$ws.Range("M21").ClearContents()  # was -265
$ws.Range("N21").Value = -2441.6666  # was None
$ws.Range("H22").Value = 536.625  # was 558.625
$ws.Range("I22").Value = 446.55554  # was 466.9
$ws.Range("J22").Value = 652.4286  # was 711.5
$ws.Range("K22").Value = 446.55554  # was 466.9
$ws.Range("L22").Value = 652.4286  # was 711.5
$ws.Range("M22").Value = -96.55554000000001  # was -116.9
$ws.Range("N22").Value = -1352.4286  # was -1411.5
$ws.Range("H86").Value = 11216.429  # was 11930.714
$ws.Range("I86").Value = 9626.75  # was 10876.75
$ws.Range("K86").Value = 9626.75  # was 10876.75
$ws.Range("M86").Value = -8503.75  # was -9753.75
$ws.Range("H89").Value = 11216.429  # was 11930.714
$ws.Range("I89").Value = 9626.75  # was 10876.75
$ws.Range("K89").Value = 48133.75  # was 54383.75
$ws.Range("M89").Value = -42517.75  # was -48767.75
$ws.Range("H105").Value = 2999.5  # was 3000
$ws.Range("I105").Value = 2999  # was 0
$ws.Range("K105").Value = 2999  # was 0
$ws.Range("M105").Value = -1252  # was None
$ws.Range("H107").Value = 2417  # was 2580.6155
$ws.Range("I107").Value = 1644.3636  # was 1779.8
$ws.Range("K107").Value = 1644.3636  # was 1779.8
$ws.Range("M107").Value = 275.6364000000001  # was 140.2
$ws.Range("H129").Value = 0  # was 40997
$ws.Range("J129").Value = 0  # was 40997
$ws.Range("L129").Value = 0  # was 40997
$ws.Range("N129").ClearContents()  # was -50997

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7001  # was 5961.6
$ws.Range("I39").Value = 7001  # was 2528.6667
$ws.Range("J39").Value = 0  # was 11111
$ws.Range("K39").Value = 21003  # was 7586.000100000001
$ws.Range("L39").Value = 0  # was 33333
$ws.Range("M39").Value = -20709  # was -7292.000100000001
$ws.Range("N39").ClearContents()  # was -33921
$ws.Range("H55").Value = 930.53845  # was 1067.3636
$ws.Range("J55").Value = 1329.125  # was 1712.8334
$ws.Range("L55").Value = 3987.375  # was 5138.5002
$ws.Range("N55").Value = -4341.375  # was -5492.5002
$ws.Range("H57").Value = 1999.5  # was 2000
$ws.Range("I57").Value = 1999.5  # was 2000
$ws.Range("K57").Value = 5998.5  # was 6000
$ws.Range("M57").Value = -5439.5  # was -5441
$ws.Range("H131").Value = 7118.2915  # was 6869.6
$ws.Range("I131").Value = 788  # was 806.8333
$ws.Range("K131").Value = 2364  # was 2420.4999
$ws.Range("M131").Value = 2676  # was 2619.5001
$ws.Range("H137").Value = 2790.4  # was 2847.4666
$ws.Range("I137").Value = 1702.875  # was 1679.1111
$ws.Range("J137").Value = 4033.2856  # was 4600
$ws.Range("K137").Value = 5108.625  # was 5037.3333
$ws.Range("L137").Value = 12099.8568  # was 13800
$ws.Range("M137").Value = -8.625  # was 62.66669999999976
$ws.Range("N137").Value = -22299.8568  # was -24000

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 55000  # was 40000
$ws.Range("H70").Value = 27839.6  # was 24158
$ws.Range("J70").Value = 22874  # was 17166
$ws.Range("L70").Value = 22874  # was 17166
$ws.Range("N70").Value = -23414  # was -17706
$ws.Range("H73").Value = 27839.6  # was 24158
$ws.Range("J73").Value = 22874  # was 17166
$ws.Range("L73").Value = 22874  # was 17166
$ws.Range("N73").Value = -24746  # was -19038
$ws.Range("H80").Value = 2941  # was 2957.8572
$ws.Range("H83").Value = 2941  # was 2957.8572
$ws.Range("H102").Value = 0  # was 2349.2
$ws.Range("I102").Value = 0  # was 1949.5
$ws.Range("J102").Value = 0  # was 2615.6667
$ws.Range("K102").Value = 0  # was 1949.5
$ws.Range("L102").Value = 0  # was 2615.6667
$ws.Range("M102").ClearContents()  # was -327.5
$ws.Range("N102").ClearContents()  # was -5859.6667
$ws.Range("H122").Value = 7278  # was 8754.637000000001
$ws.Range("I122").Value = 8172.5454  # was 9589.111000000001
$ws.Range("J122").Value = 3998  # was 4999.5
$ws.Range("K122").Value = 24517.6362  # was 28767.333
$ws.Range("L122").Value = 11994  # was 14998.5
$ws.Range("M122").Value = -22067.6362  # was -26317.333
$ws.Range("N122").Value = -16894  # was -19898.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2476.422  # was 2446.1738
$ws.Range("I22").Value = 1579.3572  # was 1562.3103
$ws.Range("K22").Value = 1579.3572  # was 1562.3103
$ws.Range("M22").Value = -1284.3572  # was -1267.3103
$ws.Range("H27").Value = 2476.422  # was 2446.1738
$ws.Range("I27").Value = 1579.3572  # was 1562.3103
$ws.Range("K27").Value = 1579.3572  # was 1562.3103
$ws.Range("M27").Value = -1472.3572  # was -1455.3103
$ws.Range("H46").Value = 4162.375  # was 3866.5557
$ws.Range("J46").Value = 4776  # was 4324.1377
$ws.Range("L46").Value = 4776  # was 4324.1377
$ws.Range("N46").Value = -5152  # was -4700.1377
$ws.Range("H47").Value = 23300  # was 25000
$ws.Range("J47").Value = 23300  # was 25000
$ws.Range("L47").Value = 23300  # was 25000
$ws.Range("N47").Value = -24280  # was -25980
$ws.Range("H52").Value = 23300  # was 25000
$ws.Range("J52").Value = 23300  # was 25000
$ws.Range("L52").Value = 23300  # was 25000
$ws.Range("N52").Value = -23766  # was -25466
$ws.Range("H63").Value = 55999.332  # was 53999.5
$ws.Range("H66").Value = 55999.332  # was 53999.5
$ws.Range("H122").Value = 6943.9  # was 6578.6665
$ws.Range("I122").Value = 5500  # was 4667.3335
$ws.Range("J122").Value = 7304.875  # was 7215.778
$ws.Range("K122").Value = 16500  # was 14002.0005
$ws.Range("L122").Value = 21914.625  # was 21647.334
$ws.Range("M122").Value = -14050  # was -11552.0005
$ws.Range("N122").Value = -26814.625  # was -26547.334
$ws.Range("H132").Value = 7058  # was 7484.696
$ws.Range("I132").Value = 6419.737  # was 6968.0625
$ws.Range("J132").Value = 8405.444  # was 8665.571
$ws.Range("K132").Value = 19259.211  # was 20904.1875
$ws.Range("L132").Value = 25216.332  # was 25996.713
$ws.Range("M132").Value = -16729.211  # was -18374.1875
$ws.Range("N132").Value = -30276.332  # was -31056.713
$ws.Range("H136").Value = 9000.799999999999  # was 11666.5
$ws.Range("I136").Value = 7500.6  # was 11249.5
$ws.Range("J136").Value = 10501  # was 11875
$ws.Range("K136").Value = 22501.8  # was 33748.5
$ws.Range("L136").Value = 31503  # was 35625
$ws.Range("M136").Value = -19951.8  # was -31198.5
$ws.Range("N136").Value = -36603  # was -40725

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 90078.5  # was 78497.5
$ws.Range("J121").Value = 90078.5  # was 78497.5
$ws.Range("L121").Value = 90078.5  # was 78497.5
$ws.Range("N121").Value = -93572.5  # was -81991.5
$ws.Range("H122").Value = 4939.067  # was 4415.1055
$ws.Range("I122").Value = 5172.8335  # was 4524.8667
$ws.Range("J122").Value = 4004  # was 4003.5
$ws.Range("K122").Value = 15518.5005  # was 13574.6001
$ws.Range("L122").Value = 12012  # was 12010.5
$ws.Range("M122").Value = -13068.5005  # was -11124.6001
$ws.Range("N122").Value = -16912  # was -16910.5
$ws.Range("H132").Value = 3940.7795  # was 3998.3965
$ws.Range("I132").Value = 2328.25  # was 2368.465
$ws.Range("K132").Value = 6984.75  # was 7105.395
$ws.Range("M132").Value = -4454.75  # was -4575.395
